$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])"
$ws.Range("B2").Value = 0.5598018648018648
$ws.Range("C2").Value = "{'selector': None, 'scaler': MinMaxScaler(), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}"
$ws.Range("D2").Value = 0.7231437256025264
$ws.Range("E2").Value = 0.5598018648018648
$ws.Range("F2").Value = 0.787878787878788
$ws.Range("G2").Value = 0.8092236053389532
$ws.Range("H2").Value = 0.6703174603174603
$ws.Range("I2").Value = 0.7647058823529411
$ws.Range("J2").Value = 0.6595744680851063
$ws.Range("K2").Value = 0.5166666666666667
$ws.Range("L2").Value = 0.8125
$ws.Range("N2").Value = "[1 1 1 1 1 1 0 1 1 0 0 1 0 1 1 1 1 1 0 1 1 0 1 0]"

$ws.Range("A3").Value = "Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])"
$ws.Range("B3").Value = 0.6000000000000001
$ws.Range("C3").Value = "{'selector': None, 'scaler': MinMaxScaler(), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}"
$ws.Range("D3").Value = 0.7526013652625287
$ws.Range("E3").Value = 0.6000000000000001
$ws.Range("F3").Value = 0.6666666666666667
$ws.Range("G3").Value = 0.7822783852303303
$ws.Range("H3").Value = 0.7138888888888889
$ws.Range("I3").Value = 0.6470588235294118
$ws.Range("J3").Value = 0.7276595744680852
$ws.Range("K3").Value = 0.55
$ws.Range("L3").Value = 0.6875
$ws.Range("N3").Value = "[0 1 0 1 1 1 1 1 1 1 0 1 1 0 1 1 1 1 0 1 1 0 0 1]"

$ws.Range("A4").Value = "Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])"
$ws.Range("B4").Value = 0.4892057942057942
$ws.Range("C4").Value = "{'selector': None, 'scaler': MinMaxScaler(), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}"
$ws.Range("D4").Value = 0.7632968482206695
$ws.Range("E4").Value = 0.4892057942057942
$ws.Range("F4").Value = 0.6060606060606061
$ws.Range("G4").Value = 0.7643994031906569
$ws.Range("H4").Value = 0.5205952380952381
$ws.Range("I4").Value = 0.7142857142857143
$ws.Range("J4").Value = 0.7755555555555556
$ws.Range("K4").Value = 0.5
$ws.Range("L4").Value = 0.5263157894736842
$ws.Range("N4").Value = "[1 1 1 1 0 0 0 1 1 1 1 0 0 0 1 0 0 0 1 0 1 1 1 1]"

$ws.Range("A5").Value = "Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])"
$ws.Range("B5").Value = 0.5007209457209457
$ws.Range("C5").Value = "{'selector': None, 'scaler': MinMaxScaler(), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}"
$ws.Range("D5").Value = 0.7321640629636167
$ws.Range("E5").Value = 0.5007209457209457
$ws.Range("F5").Value = 0.6923076923076924
$ws.Range("G5").Value = 0.8036962325474752
$ws.Range("H5").Value = 0.5795238095238096
$ws.Range("I5").Value = 0.75
$ws.Range("J5").Value = 0.6857142857142857
$ws.Range("K5").Value = 0.4666666666666667
$ws.Range("L5").Value = 0.6428571428571429
$ws.Range("N5").Value = "[0 1 1 0 0 0 0 1 1 0 1 0 1 1 0 0 1 0 1 0 1 1 0 1]"

$ws.Range("A6").Value = "Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])"
$ws.Range("B6").Value = 0.552121212121212
$ws.Range("C6").Value = "{'selector': None, 'scaler': MinMaxScaler(), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}"
$ws.Range("D6").Value = 0.7727446463776266
$ws.Range("E6").Value = 0.552121212121212
$ws.Range("F6").Value = 0.5833333333333334
$ws.Range("G6").Value = 0.8477183787689755
$ws.Range("H6").Value = 0.6248809523809524
$ws.Range("I6").Value = 0.5384615384615384
$ws.Range("J6").Value = 0.7134615384615384
$ws.Range("K6").Value = 0.5499999999999999
$ws.Range("L6").Value = 0.6363636363636364
$ws.Range("N6").Value = "[1 0 1 0 1 0 1 0 0 0 0 1 0 1 1 0 0 0 1 1 1 1 1 1]"
